$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 5
$ws.Range("F4").Value = -1
$ws.Range("F12").Value = -9
$ws.Range("F22").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("F28").Value = -3
$ws.Range("F29").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("F41").Value = -1
$ws.Range("F42").Value = 0
$ws.Range("F45").Value = 1
$ws.Range("F47").Value = 2
$ws.Range("F50").Value = 1
$ws.Range("F56").Value = 0
$ws.Range("F58").Value = 1
$ws.Range("F59").Value = -3
$ws.Range("F60").Value = 0
$ws.Range("F61").Value = 3
$ws.Range("F63").Value = -1
$ws.Range("F64").Value = -2
$ws.Range("F66").Value = -1
